$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 5 (shifts 신한스팩12호 and below down by one)
$ws.Rows.Item(5).Insert()

# Fill the new row 5 with the new SPAC offering data
$ws.Cells.Item(5, 1).Value = "신한스팩13호"
$ws.Cells.Item(5, 2).Value = "2024.04.04~04.05"
$ws.Cells.Item(5, 3).Value = "2,000~2,000"
$ws.Cells.Item(5, 4).Value = "-"
$ws.Cells.Item(5, 5).Value = 6000
$ws.Cells.Item(5, 6).Value = "신한투자증권"

# Delete the last row (was row 21, now row 22 after the insert) -- the 이에이트 row
$ws.Rows.Item(22).Delete()
